$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from row 2 to row 33
# with the new date serial value 46062 (2026-02-09), replacing 46061.
$ws.Range("C2:C33").Value = 46062
